$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity

$rows = @(
    @("FAPs","Fgf17","Fgfr1","ECs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,5.39270102941,48.53430926468999,0.0202329253530045,0.0202329253530045),
    @("FAPs","Fgf17","Fgfr1","FAPs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,118.0082709195989,1062.07443827639,0.4427563337058576,0.4427563337058578),
    @("FAPs","Fgf17","Fgfr1","M2",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,0.8178189009914444,7.360370108922999,0.003068382371986681,0.003068382371986682),
    @("FAPs","Fgf17","Fgfr1","sCs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,76.50370833413376,688.5333750072039,0.2870349777432233,0.2870349777432233),
    @("sCs","Fgf17","Fgfr1","ECs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,1.76803975082,15.91235775738,0.006633524852276025,0.006633524852276028),
    @("sCs","Fgf17","Fgfr1","FAPs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,38.68994642453112,348.20951782078,0.1451611713035993,0.1451611713035994),
    @("sCs","Fgf17","Fgfr1","M2",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,0.2681284050495556,2.413155645446,0.00100599346687338,0.001005993466873381),
    @("sCs","Fgf17","Fgfr1","sCs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,25.08234680213422,225.741121219208,0.09410669120317908,0.09410669120317912)
)

$rowIndex = 2
foreach ($rowData in $rows) {
    $colIndex = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
